# CambioTasa.xlsx update ("update entregable 1, 2 y 3")
#
# The sheet header cell D1 was relabelled from "fijo" to "tasa" (the
# underlying data/account rows are unchanged). Also nudge the saved
# selection/active cell to F8 to match the author's last selection when
# the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "tasa"

$ws.Range("F8").Select()
